$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5","D6","D8","D14","D19","D20","D21","D22","D23","D24","D25","D26","D28","D29","D31","D32","D34","D35","D38","D39","D40","D42","D43","D45","D46","D48","D49","D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.247.00"
$ws.Range("E2").Value = "  +0.65%  "

$ws.Range("D3").Value = "3.477.82"
$ws.Range("E3").Value = "  -0.34%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "593.66"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").Value = "178.72"
$ws.Range("E6").Value = "  +4.18%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "0.597"
$ws.Range("E8").Value = "  +1.16%  "

$ws.Range("D9").Value = "3.476.33"
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("E10").Value = "  +4.04%  "

$ws.Range("E11").Value = "  -2.25%  "

$ws.Range("E12").Value = "  +0.97%  "

$ws.Range("D13").Value = "4.076.52"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").Value = "31.97"
$ws.Range("E14").Value = "  +10.84%  "

$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").Value = "67.270.14"
$ws.Range("E16").Value = "  +0.74%  "

$ws.Range("E17").Value = "  -0.84%  "

$ws.Range("D18").Value = "3.474.55"
$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").Value = "6.25"
$ws.Range("E19").Value = "  -0.17%  "

$ws.Range("D20").Value = "14.29"
$ws.Range("E20").Value = "  +2.06%  "

$ws.Range("D21").Value = "388.33"
$ws.Range("E21").Value = "  -0.99%  "

$ws.Range("D22").Value = "7.95"
$ws.Range("E22").Value = "  +0.25%  "

$ws.Range("D23").Value = "74.15"
$ws.Range("E23").Value = "  +1.83%  "

$ws.Range("D24").Value = "0.996"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("B25").Value = "LEO"
$ws.Range("C25").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D25").Value = "5.72"
$ws.Range("E25").Value = "  +0.49%  "

$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.535"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("E27").Value = "  +0.12%  "

$ws.Range("D28").Value = "10.38"
$ws.Range("E28").Value = "  +2.38%  "

$ws.Range("D29").Value = "0.174"
$ws.Range("E29").Value = "  -3.41%  "

$ws.Range("E30").Value = "  +0.24%  "

$ws.Range("D31").Value = "6.14"
$ws.Range("E31").Value = "  -0.34%  "

$ws.Range("D32").Value = "1.41"
$ws.Range("E32").Value = "  -1.13%  "

$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").Value = "23.51"
$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("D35").Value = "7.35"
$ws.Range("E35").Value = "  +0.26%  "

$ws.Range("E36").Value = "  +0.04%  "

$ws.Range("E37").Value = "  -1.69%  "

$ws.Range("D38").Value = "163.79"
$ws.Range("E38").Value = "  +0.66%  "

$ws.Range("D39").Value = "0.869"
$ws.Range("E39").Value = "  -0.83%  "

$ws.Range("D40").Value = "1.87"
$ws.Range("E40").Value = "  -1.20%  "

$ws.Range("E41").Value = "  +6.52%  "

$ws.Range("D42").Value = "6.83"
$ws.Range("E42").Value = "  -1.28%  "

$ws.Range("D43").Value = "4.62"
$ws.Range("E43").Value = "  -0.50%  "

$ws.Range("D44").Value = "2.835.54"
$ws.Range("E44").Value = "  +0.89%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "26.98"
$ws.Range("E45").Value = "  -0.70%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "26.15"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("E47").Value = "  -2.51%  "

$ws.Range("D48").Value = "41.57"
$ws.Range("E48").Value = "  -2.79%  "

$ws.Range("D49").Value = "0.0300"
$ws.Range("E49").Value = "  -0.41%  "

$ws.Range("D50").Value = "337.29"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("E51").Value = "  -2.69%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}